# Evaluation Table.xlsx — add new dataset evaluation block + classify-label column
# (commit: "Update Dataset, Add classify function")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Pass 1: create the new shared strings in the same order they were first
# authored, so the shared-strings table lines up index-for-index with the
# target workbook (20..29).
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = "new Pyfeat Data 사용"
$ws.Range("A17").Value = "train: test = 7:3"
$ws.Range("A18").Value = "train: test = 7:3 -kfold"
$ws.Range("A19").Value = "train: test = 8:2"
$ws.Range("A20").Value = "train: test = 8:2 -kfold"
$ws.Range("A21").Value = "train: test = 5:5"
$ws.Range("I22").Value = " "
$ws.Range("I21").Value = "정확도 66 ~72사이에서 변동"
$ws.Range("A22").Value = "train: test = 5:5 -kfold"
$ws.Range("L2").Value = "레이어 설정 base_model"

# ---------------------------------------------------------------------------
# Pass 2: new header row (row 16) — reuses the existing metric-label strings.
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = "accuracy"
$ws.Range("C16").Value = "sensitivity"
$ws.Range("D16").Value = "specificity"
$ws.Range("E16").Value = "PPV"
$ws.Range("F16").Value = "NPV"
$ws.Range("G16").Value = "MCC"
$ws.Range("H16").Value = "AUC"

# ---------------------------------------------------------------------------
# Pass 3: numeric evaluation data for the new dataset block (rows 17-22).
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = 0.696
$ws.Range("C17").Value = 0.833
$ws.Range("D17").Value = 0.5333
$ws.Range("E17").Value = 0.681
$ws.Range("F17").Value = 0.72
$ws.Range("G17").Value = 0.387
$ws.Range("H17").Value = 0.725

$ws.Range("B18").Value = 0.972
$ws.Range("C18").Value = 0.975
$ws.Range("D18").Value = 0.966
$ws.Range("E18").Value = 0.981
$ws.Range("F18").Value = 0.957
$ws.Range("G18").Value = 0.94
$ws.Range("H18").Value = 0.98

$ws.Range("B19").Value = 0.679
$ws.Range("C19").Value = 0.88
$ws.Range("D19").Value = 0.5
$ws.Range("E19").Value = 0.611
$ws.Range("F19").Value = 0.823
$ws.Range("G19").Value = 0.406
$ws.Range("H19").Value = 0.745

$ws.Range("B20").Value = 0.994
$ws.Range("C20").Value = 0.993
$ws.Range("D20").Value = 0.1
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.983
$ws.Range("G20").Value = 0.988
$ws.Range("H20").Value = 1

$ws.Range("B21").Value = 0.722
$ws.Range("C21").Value = 0.854
$ws.Range("D21").Value = 0.593
$ws.Range("E21").Value = 0.67
$ws.Range("F21").Value = 0.88
$ws.Range("G21").Value = 0.463
$ws.Range("H21").Value = 0.801

$ws.Range("B22").Value = 0.983
$ws.Range("C22").Value = 0.985
$ws.Range("D22").Value = 0.985
$ws.Range("E22").Value = 0.983
$ws.Range("F22").Value = 0.983
$ws.Range("G22").Value = 0.969
$ws.Range("H22").Value = 0.997

# ---------------------------------------------------------------------------
# Pass 4: classify-label column L — tags every new row (and the original
# CNN_base_model_k_fold row 2) with the "레이어 설정 base_model" label.
# ---------------------------------------------------------------------------
$ws.Range("L17").Value = "레이어 설정 base_model"
$ws.Range("L18").Value = "레이어 설정 base_model"
$ws.Range("L19").Value = "레이어 설정 base_model"
$ws.Range("L20").Value = "레이어 설정 base_model"
$ws.Range("L21").Value = "레이어 설정 base_model"
$ws.Range("L22").Value = "레이어 설정 base_model"

# ---------------------------------------------------------------------------
# Final selection, matching the saved cursor position in the target file.
# ---------------------------------------------------------------------------
$ws.Range("F21").Select()

Write-Output "edit applied"
